$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# League base update (08-05-2024 20:15): two pairs of match rows were
# re-ordered (216<->217 and 219<->220) and several odds columns were
# refreshed on rows 216-222.
# ---------------------------------------------------------------------------

# Swap the full content (id, teams, odds) of rows 216 and 217, cell by cell,
# using a scratch cell far outside the used range as temporary storage so
# that cell types (text vs number) and styles are preserved exactly.
$swapCols = @("B","E","F","J","K","L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $swapCols) {
    $ws.Range($col + "216").Copy($ws.Range("ZZ1"))
    $ws.Range($col + "217").Copy($ws.Range($col + "216"))
    $ws.Range("ZZ1").Copy($ws.Range($col + "217"))
}
$ws.Range("ZZ1").ClearContents()

# Swap the full content of rows 219 and 220 the same way.
foreach ($col in $swapCols) {
    $ws.Range($col + "219").Copy($ws.Range("ZZ1"))
    $ws.Range($col + "220").Copy($ws.Range($col + "219"))
    $ws.Range("ZZ1").Copy($ws.Range($col + "220"))
}
$ws.Range("ZZ1").ClearContents()

# Refresh odds lines that moved along with (but not purely mirrored by) the
# row swap above.
$ws.Range("Q216").Value = 1.98
$ws.Range("R216").Value = 1.92

$ws.Range("Q217").Value = 1.84
$ws.Range("R217").Value = 2.06

$ws.Range("T219").Value = 1.925
$ws.Range("U219").Value = 1.925

# Row 218 odds refresh.
$ws.Range("P218").Value = 0
$ws.Range("Q218").Value = 1.8
$ws.Range("R218").Value = 2.1
$ws.Range("T218").Value = 1.95
$ws.Range("U218").Value = 1.9

# Row 221 odds refresh.
$ws.Range("Q221").Value = 1.89
$ws.Range("R221").Value = 2.01

# Row 222 odds refresh.
$ws.Range("Q222").Value = 1.99
$ws.Range("R222").Value = 1.91
$ws.Range("T222").Value = 1.85
$ws.Range("U222").Value = 2
